# Weekly driver report update for 2025-04-29
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments ---
# Raw stored width = ColumnWidth + 0.8333333333333334, so subtract that offset
# to land exactly on the desired stored widths.
$offset = 0.8333333333333334
$ws.Columns.Item(2).ColumnWidth  = (14 - $offset)   # B: 15 -> 14
$ws.Columns.Item(5).ColumnWidth  = (14 - $offset)   # E: 16 -> 14
$ws.Columns.Item(6).ColumnWidth  = (11 - $offset)   # F: 2  -> 11
$ws.Columns.Item(7).ColumnWidth  = (43 - $offset)   # G: 2  -> 43
$ws.Columns.Item(8).ColumnWidth  = (14 - $offset)   # H: 2  -> 14
$ws.Columns.Item(9).ColumnWidth  = (30 - $offset)   # I: 2  -> 30
$ws.Columns.Item(10).ColumnWidth = (16 - $offset)   # J: 2  -> 16

# --- Updated "Bad Drivers" figures ---
$ws.Range("C3").Value = 2385
$ws.Range("D3").Value = 33.2
$ws.Range("C4").Value = 2385

# --- Rebuild the "Good Drivers" header row (row 11) ---
# New layout spans A:J and drops the previous bold/border formatting.
$ws.Range("A11:J11").ClearFormats()
$ws.Range("A11").Value = "adapter-driver"
$ws.Range("B11").Value = "good sum"
$ws.Range("C11").Value = "critical sum"
$ws.Range("D11").Value = "warning sum"
$ws.Range("E11").Value = "client count"
$ws.Range("F11").Value = "total sum"
$ws.Range("G11").Value = "adapter"
$ws.Range("H11").Value = "driver"
$ws.Range("I11").Value = "good roaming calculation (%)"
$ws.Range("J11").Value = "driver vintage"

# --- Rebuild the "Good Drivers" data row (row 12) ---
$ws.Range("A12:J12").ClearFormats()
$ws.Range("A12").Value = "MediaTek Wi-Fi 6 MT7921 Wireless LAN Card - 22.30.1.1339"
$ws.Range("B12").Value = 10921
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 35
$ws.Range("F12").Value = 10923
$ws.Range("G12").Value = "mediatek wi-fi 6 mt7921 wireless lan card"
$ws.Range("H12").Value = "22.30.1.1339"
$ws.Range("I12").Value = 100

# Force J12 to remain a plain text value instead of being auto-parsed as a date:
# enter it as a formula returning a string (bypasses literal-entry date detection
# and avoids allocating a brand new number-format style), then convert the
# formula result down to a plain value via copy / paste-values.
$ws.Range("J12").Formula = '="2022-08-18"'
$ws.Range("J12").Copy() | Out-Null
$ws.Range("J12").PasteSpecial(-4163) | Out-Null  # xlPasteValues
$excel.CutCopyMode = $false

# Rows 13-17 are already empty, so clearing them ensures no stray content remains
# and lets the sheet dimension shrink back down to J12 automatically.
$ws.Range("A13:J17").Clear()
